$d = $word.ActiveDocument

# The final section of the document (a Jekyll-site export) ends with a
# bibliography entry followed by a blank paragraph, a "Ver no Jupiter..."
# line, and a "(c) 2020 ..." footer line. This edit strips the trailing
# blank paragraph plus those two footer paragraphs, leaving the
# bibliography entry directly followed by the original blank spacer
# paragraph and the final page-break paragraph.

$target = "Ver no Jupiter Salvar em pdf Salvar em docx"

$finder = $d.Content
$found = $finder.Find.Execute($target, $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'Ver no Jupiter...' paragraph"
}
$targetStart = $finder.Start

# Resolve the paragraph index of the matched range (use a strict upper
# bound so the shared boundary point between consecutive paragraphs
# doesn't cause the preceding paragraph to match instead).
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $targetStart -and $p.Range.End -gt $targetStart) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq 0) {
    throw "Could not resolve paragraph index for the matched range"
}

# Remove the blank paragraph right before it, the "Ver no Jupiter..."
# paragraph itself, and the copyright paragraph right after it.
$firstPara = $d.Paragraphs.Item($targetIndex - 1)
$lastPara = $d.Paragraphs.Item($targetIndex + 1)

$deleteRange = $d.Range($firstPara.Range.Start, $lastPara.Range.End)
$deleteRange.Delete()
